# TestCaseToExecute.xlsx update:
# Enable (activate) the run flag for TestCaseID 2 (row 3) so that the
# login module / Admin - user management test case is included in the run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# RunFlag (column B) for TestCaseID 2 (row 3): 0 -> 1 (active)
$ws.Range("B3").Value = 1

# Leave the cursor/selection on the cell that was just edited, as in the
# authored workbook.
$ws.Range("B3").Select()
